$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the pre-existing blank F1 cell (engine round-trips an empty
# shared-string cell with no <v> into string index 0 unless we touch it
# explicitly) - clearing keeps it truly empty, matching the source file.
$ws.Range("F1").ClearContents()

# Add "NA" in the duplicate_image_filename column (E) for every row of
# the first data table (rows 2-21).
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
